$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'products__item\ in-stock\ products__item_4-in-row"]:nth-child(3) [type="button',
    'Hot! New! Best choiceiPhone',
    'Increased',
    'Buy',
    '\31 52175-case-658',
    '\31 52176-case-660',
    '512 GB',
    'Natural Titanium'
)

foreach ($r in 4, 5) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($r, $i + 1).Value = $values[$i]
    }
}
